# Scheduled market-data refresh: update computed price/profit columns (H:N)
# across the Leve profitability sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 45
$ws.Range("H45").Value = 1285.6
$ws.Range("I45").Value = 1285.6
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 3856.8
$ws.Range("L45").Value = 0
$ws.Range("M45").Value = -3664.8
$ws.Range("N45").Value = $null

# Row 112
$ws.Range("H112").Value = 92972.82000000001
$ws.Range("J112").Value = 2041.6
$ws.Range("L112").Value = 6124.799999999999
$ws.Range("N112").Value = -8340.799999999999

# Row 113
$ws.Range("H113").Value = 76927440
$ws.Range("J113").Value = 5071.5713
$ws.Range("L113").Value = 5071.5713
$ws.Range("N113").Value = -11579.5713

# Row 137
$ws.Range("H137").Value = 1825.4546
$ws.Range("I137").Value = 1658
$ws.Range("K137").Value = 4974
$ws.Range("M137").Value = -2424

# Row 138
$ws.Range("H138").Value = 2939.0806
$ws.Range("I138").Value = 1166.7391
$ws.Range("K138").Value = 3500.2173
$ws.Range("M138").Value = 1639.7827

$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 240
$ws.Range("I4").Value = 110
$ws.Range("J4").Value = 500
$ws.Range("K4").Value = 110
$ws.Range("L4").Value = 500
$ws.Range("M4").Value = 6
$ws.Range("N4").Value = -732

# Row 32
$ws.Range("H32").Value = 1830.5834
$ws.Range("I32").Value = 1830.5834
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 1830.5834
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -1543.5834
$ws.Range("N32").Value = $null

# Row 45
$ws.Range("H45").Value = 2490.818
$ws.Range("I45").Value = 2166.6667
$ws.Range("J45").Value = 2879.8
$ws.Range("K45").Value = 2166.6667
$ws.Range("L45").Value = 2879.8
$ws.Range("M45").Value = -1789.6667
$ws.Range("N45").Value = -3633.8

# Row 74
$ws.Range("H74").Value = 3868.4
$ws.Range("I74").Value = 2812.7
$ws.Range("K74").Value = 2812.7
$ws.Range("M74").Value = -1938.7

# Row 77
$ws.Range("H77").Value = 3868.4
$ws.Range("I77").Value = 2812.7
$ws.Range("K77").Value = 14063.5
$ws.Range("M77").Value = -9695.5

# Row 102
$ws.Range("H102").Value = 57761.07
$ws.Range("I102").Value = 819.6667
$ws.Range("J102").Value = 160255.6
$ws.Range("K102").Value = 819.6667
$ws.Range("L102").Value = 160255.6
$ws.Range("M102").Value = 802.3333
$ws.Range("N102").Value = -163499.6

# Row 110
$ws.Range("H110").Value = 85815.336
$ws.Range("I110").Value = 52628.4
$ws.Range("K110").Value = 52628.4
$ws.Range("M110").Value = -50583.4

# Row 122
$ws.Range("H122").Value = 27781160
$ws.Range("I122").Value = 33336392
$ws.Range("K122").Value = 100009176
$ws.Range("M122").Value = -100006726

# Row 132
$ws.Range("H132").Value = 43480976
$ws.Range("I132").Value = 55557976
$ws.Range("K132").Value = 166673928
$ws.Range("M132").Value = -166671398

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 28665.666
$ws.Range("I75").Value = 27998.5
$ws.Range("K75").Value = 27998.5
$ws.Range("M75").Value = -27062.5

# Row 78
$ws.Range("H78").Value = 28665.666
$ws.Range("I78").Value = 27998.5
$ws.Range("K78").Value = 83995.5
$ws.Range("M78").Value = -79315.5

# Row 134
$ws.Range("H134").Value = 2013.3784
$ws.Range("I134").Value = 1716.0646
$ws.Range("K134").Value = 5148.1938
$ws.Range("M134").Value = -2613.1938

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2860.9553
$ws.Range("I31").Value = 2217.64
$ws.Range("J31").Value = 3243.8809
$ws.Range("K31").Value = 2217.64
$ws.Range("L31").Value = 3243.8809
$ws.Range("M31").Value = -1922.64
$ws.Range("N31").Value = -3833.8809

# Row 34
$ws.Range("H34").Value = 2860.9553
$ws.Range("I34").Value = 2217.64
$ws.Range("J34").Value = 3243.8809
$ws.Range("K34").Value = 2217.64
$ws.Range("L34").Value = 3243.8809
$ws.Range("M34").Value = -2015.64
$ws.Range("N34").Value = -3647.8809

# Row 122
$ws.Range("H122").Value = 2670.2856
$ws.Range("I122").Value = 2262.75
$ws.Range("J122").Value = 3213.6667
$ws.Range("K122").Value = 6788.25
$ws.Range("L122").Value = 9641.000100000001
$ws.Range("M122").Value = -4338.25
$ws.Range("N122").Value = -14541.0001

# Row 132
$ws.Range("H132").Value = 3154.6667
$ws.Range("I132").Value = 3067.8125
$ws.Range("K132").Value = 9203.4375
$ws.Range("M132").Value = -6673.4375

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 388.64
$ws.Range("I122").Value = 380.8889
$ws.Range("K122").Value = 3428.0001
$ws.Range("M122").Value = -978.0000999999997

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 4433.875
$ws.Range("I80").Value = 4196.8
$ws.Range("K80").Value = 4196.8
$ws.Range("M80").Value = -3198.8

# Row 83
$ws.Range("H83").Value = 4433.875
$ws.Range("I83").Value = 4196.8
$ws.Range("K83").Value = 20984
$ws.Range("M83").Value = -15992

# Row 132
$ws.Range("H132").Value = 3967.0938
$ws.Range("I132").Value = 3258.739
$ws.Range("J132").Value = 5777.3335
$ws.Range("K132").Value = 9776.217000000001
$ws.Range("L132").Value = 17332.0005
$ws.Range("M132").Value = -7246.217000000001
$ws.Range("N132").Value = -22392.0005

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 9944.556
$ws.Range("I61").Value = 8917.368
$ws.Range("J61").Value = 12384.125
$ws.Range("K61").Value = 8917.368
$ws.Range("L61").Value = 12384.125
$ws.Range("M61").Value = -8715.368
$ws.Range("N61").Value = -12788.125

# Row 63
$ws.Range("H63").Value = 45831.668
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 45831.668
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 45831.668
$ws.Range("M63").Value = $null
$ws.Range("N63").Value = -47329.668

# Row 66
$ws.Range("H66").Value = 45831.668
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 45831.668
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 137495.004
$ws.Range("M66").Value = $null
$ws.Range("N66").Value = -144983.004

# Row 74
$ws.Range("H74").Value = 51364.832
$ws.Range("I74").Value = 50038.8
$ws.Range("J74").Value = 57995
$ws.Range("K74").Value = 50038.8
$ws.Range("L74").Value = 57995
$ws.Range("M74").Value = -49040.8
$ws.Range("N74").Value = -59991

# Row 77
$ws.Range("H77").Value = 51364.832
$ws.Range("I77").Value = 50038.8
$ws.Range("J77").Value = 57995
$ws.Range("K77").Value = 150116.4
$ws.Range("L77").Value = 173985
$ws.Range("M77").Value = -145124.4
$ws.Range("N77").Value = -183969

# Row 94
$ws.Range("H94").Value = 50017476
$ws.Range("J94").Value = 50017476
$ws.Range("L94").Value = 50017476
$ws.Range("N94").Value = -50018828

# Row 113
$ws.Range("H113").Value = 9944.556
$ws.Range("I113").Value = 8917.368
$ws.Range("J113").Value = 12384.125
$ws.Range("K113").Value = 8917.368
$ws.Range("L113").Value = 12384.125
$ws.Range("M113").Value = -6747.368
$ws.Range("N113").Value = -16724.125

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1955.5454
$ws.Range("I126").Value = 1693.8572
$ws.Range("K126").Value = 5081.571599999999
$ws.Range("M126").Value = -2611.571599999999

# Row 132
$ws.Range("H132").Value = 3009.9424
$ws.Range("I132").Value = 2667.5
$ws.Range("J132").Value = 4893.375
$ws.Range("K132").Value = 8002.5
$ws.Range("L132").Value = 14680.125
$ws.Range("M132").Value = -5472.5
$ws.Range("N132").Value = -19740.125
